$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: reorder/replace field labels ---
# old order: name, email, password, cpf, phone, cargo
# new order: name, email, phone, cpf, cargo, status
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "email"
$ws.Range("C1").Value = "phone"
$ws.Range("D1").Value = "cpf"
$ws.Range("E1").Value = "cargo"
$ws.Range("F1").Value = "status"

# Header row no longer carries the wrap-text style; revert to default
$ws.Range("A1:F1").Style = "Normal"

# Row 1 no longer has an explicit custom height - autofit back to default
$ws.Rows(1).AutoFit()

# --- Body rows: the "password" column (B) is no longer part of the sheet ---
$ws.Range("B2:B4").Clear()

# --- Selection moves ---
$ws.Range("D8").Select()
